# Applies the "Fixed Lists of Recommended Factors" edit:
#  - Column M ("Drop:") is trimmed down to a curated list of 17 factors to
#    drop (header + 17 rows, M8:M25).
#  - A brand-new column P ("Use:") is added, listing the remaining 34
#    factors to use (header + 34 rows, P8:P42).
#  - New shared strings introduced: "Use:", "bcgv", "dtpv", "fullv",
#    "mslv", "poliov".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column P ("Use:") header, written first so new shared strings land
# in the same order as the reference edit (Use: right after Drop:). -----
$ws.Cells.Item(8, 16).Value  = "Use:"
$ws.Cells.Item(8, 16).Font.Bold = $true

# ---- Column M ("Drop:") -----------------------------------------------
# Row 8 (header "Drop:") is unchanged.
# Rows 9-13 become the short vaccine-code list.
$ws.Cells.Item(9, 13).Value  = "bcgv"
$ws.Cells.Item(10, 13).Value = "dtpv"
$ws.Cells.Item(11, 13).Value = "fullv"
$ws.Cells.Item(12, 13).Value = "mslv"
$ws.Cells.Item(13, 13).Value = "poliov"
# Rows 14-22 keep their previous values (cmr, stunt, uweight, pnmr,
# NUTSTUNTINGPREV, NUTUNDERWEIGHTPREV, MEDS1_01_03, WAS_0000000001,
# EQ_HANDWASHING) - nothing to change there.
# Rows 23-25 shift down by one position in the drop list.
$ws.Cells.Item(23, 13).Value = "EQ_OPENDEFECATION "
$ws.Cells.Item(24, 13).Value = "WAS_0000000002 "
$ws.Cells.Item(25, 13).Value = "WHOSIS_000011 "

$useValues = @(
    "WHS4_543 ",
    "WHS4_100 ",
    "full ",
    "WHS4_117 ",
    "WHS4_129 ",
    "measlesv ",
    "WHS8_110 ",
    "PCV3 ",
    "WHS4_544 ",
    "ROTAC ",
    "NUTRITION_564 ",
    "MCV2 ",
    "WHS4_128 ",
    "pncall5 ",
    "pncall3 ",
    "WHOSIS_000006 ",
    "LBW_NUMBER ",
    "LBW_PREVALENCE ",
    "WHOSIS_000005 ",
    "NUTRITION_WH2 ",
    "NUTRITION_HA_2 ",
    "NUTRITION_WA_2 ",
    "NUTRITION_WH_2 ",
    "GHED_CHEGDP_SHA2011 ",
    "WHS9_85 ",
    "MDG_0000000026 ",
    "WHS_PBR ",
    "WHS9_95 ",
    "WSH_2 ",
    "WSH_3 ",
    "WSH_SANITATION_SAFELY_MANAGED ",
    "M_Est_smk_curr ",
    "M_Est_smk_daily ",
    "TOBACCO_0000000192 "
)

$row = 9
foreach ($val in $useValues) {
    $ws.Cells.Item($row, 16).Value = $val
    $row = $row + 1
}

# ---- View bookkeeping ---------------------------------------------------
$ws.Range("L35").Select()
